$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195, shifting the existing rows 195-198 down to 196-199
$ws.Rows("195:195").Insert()

# Populate the newly inserted row 195 with the new weekly price entry
$ws.Range("A195").Value = 8
$ws.Range("B195").Value = "Terminal La Palmera de La Serena"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 45239
$ws.Range("E195").Value = 4
$ws.Range("F195").Value = 100114007
$ws.Range("G195").Value = "Jengibre"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 300
$ws.Range("K195").Value = 24000
$ws.Range("L195").Value = 25000
$ws.Range("M195").Value = 24500
$ws.Range("N195").Value = "$/caja 13 kilos"
$ws.Range("O195").Value = "Perú"
$ws.Range("P195").Value = 1885
$ws.Range("Q195").Value = 13
$ws.Range("R195").Value = "Hortaliza"
